# Update master to output generated at 9a8706d
$d = $word.ActiveDocument

$replacements = @(
    @{old = "90÷2=45, 0"; new = "75÷2=37, 1"},
    @{old = "33÷7=4, 5"; new = "71÷2=35, 1"},
    @{old = "60÷8=7, 4"; new = "71÷8=8, 7"},
    @{old = "19÷8=2, 3"; new = "98÷8=12, 2"},
    @{old = "90÷9=10, 0"; new = "68÷4=17, 0"},
    @{old = "86÷6=14, 2"; new = "56÷7=8, 0"},
    @{old = "62÷8=7, 6"; new = "62÷3=20, 2"},
    @{old = "69÷4=17, 1"; new = "41÷9=4, 5"},
    @{old = "95÷8=11, 7"; new = "97÷5=19, 2"},
    @{old = "58÷2=29, 0"; new = "47÷3=15, 2"},
    @{old = "49÷5=9, 4"; new = "23÷2=11, 1"},
    @{old = "99÷8=12, 3"; new = "71÷7=10, 1"},
    @{old = "81÷5=16, 1"; new = "46÷8=5, 6"},
    @{old = "87÷8=10, 7"; new = "57÷3=19, 0"},
    @{old = "41÷4=10, 1"; new = "18÷4=4, 2"},
    @{old = "94÷5=18, 4"; new = "76÷7=10, 6"},
    @{old = "97÷3=32, 1"; new = "84÷6=14, 0"},
    @{old = "13÷6=2, 1"; new = "98÷2=49, 0"},
    @{old = "94÷6=15, 4"; new = "80÷7=11, 3"},
    @{old = "49÷9=5, 4"; new = "73÷4=18, 1"},
    @{old = "96÷6=16, 0"; new = "54÷2=27, 0"},
    @{old = "80÷3=26, 2"; new = "21÷3=7, 0"},
    @{old = "79÷6=13, 1"; new = "24÷9=2, 6"},
    @{old = "47÷4=11, 3"; new = "44÷8=5, 4"},
    @{old = "66÷6=11, 0"; new = "92÷9=10, 2"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
